{"js": "// The \"Synthetical View\" summary table gains one more results row\n// (duplicating the last row's look), reporting the \"Cross Check\n// Acronyms\" rule a second time with Total=218 / Failed=218.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\ntable.addRows(\"End\", 1, [\n  [\"Ad_T_A445792it-C.docx\", \"Cross Check Acronyms\", \"218\", \"218\"]\n]);\nawait context.sync();\n", "ps1": "# The \"Synthetical View\" summary table gains one more results row\n# (duplicating the last row's look), reporting the \"Cross Check\n# Acronyms\" rule a second time with Total=218 / Failed=218.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRow = $t.Rows.Add()\n$i = $newRow.Index\n\n$t.Cell($i, 1).Range.Text = \"Ad_T_A445792it-C.docx\"\n$t.Cell($i, 2).Range.Text = \"Cross Check Acronyms\"\n$t.Cell($i, 3).Range.Text = \"218\"\n$t.Cell($i, 4).Range.Text = \"218\"\n"}
